$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Replace the header cell A1 with the concatenated table contents text
$ws.Range("A1").Value = "Grade Subjects Contents 1 English Quiz 1 Math Quiz,Worksheet 2 Science Flashcards"

# Update the active selection on Sheet2 to D12
$ws.Activate()
$ws.Range("D12").Select()
